$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026941710569776
$ws.Range("D2").Value = 1.031819275381261
$ws.Range("E2").Value = 1.027099454915871
$ws.Range("F2").Value = 1.025466810597581
$ws.Range("I2").Value = 1.033262555753757
$ws.Range("J2").Value = 1.032102346633824
$ws.Range("K2").Value = 1.034626292176537
$ws.Range("L2").Value = 1.029920171200744
$ws.Range("M2").Value = 1.028292297631711
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027851296594316
$ws.Range("D3").Value = 1.032502739018733
$ws.Range("E3").Value = 1.027870803901045
$ws.Range("F3").Value = 1.027010032195512
$ws.Range("I3").Value = 1.033467927207493
$ws.Range("J3").Value = 1.032652111951126
$ws.Range("K3").Value = 1.035118734488251
$ws.Range("L3").Value = 1.030499272198093
$ws.Range("M3").Value = 1.029640831784896
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.0284399425353
$ws.Range("D4").Value = 1.032945006340044
$ws.Range("E4").Value = 1.028370367702346
$ws.Range("F4").Value = 1.028008677709983
$ws.Range("I4").Value = 1.033599659272663
$ws.Range("J4").Value = 1.033007315989617
$ws.Range("K4").Value = 1.035436721433211
$ws.Range("L4").Value = 1.030873781325145
$ws.Range("M4").Value = 1.030513021057768
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028687428965312
$ws.Range("D5").Value = 1.033130939751491
$ws.Range("E5").Value = 1.028580491256475
$ws.Range("F5").Value = 1.028428531665195
$ws.Range("I5").Value = 1.033654762261839
$ws.Range("J5").Value = 1.033156516508562
$ws.Range("K5").Value = 1.035570245727098
$ws.Range("L5").Value = 1.031031174835246
$ws.Range("M5").Value = 1.030879597063609
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028728984174363
$ws.Range("D6").Value = 1.033162159016621
$ws.Range("E6").Value = 1.028615778177016
$ws.Range("F6").Value = 1.028499028499897
$ws.Range("I6").Value = 1.033663998037413
$ws.Range("J6").Value = 1.033181560464866
$ws.Range("K6").Value = 1.035592655802824
$ws.Range("L6").Value = 1.031057598948595
$ws.Range("M6").Value = 1.030941141487732
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028443249384007
$ws.Range("D7").Value = 1.032947490775675
$ws.Range("E7").Value = 1.028373174964232
$ws.Range("F7").Value = 1.028014287721331
$ws.Range("I7").Value = 1.033600396650823
$ws.Range("J7").Value = 1.033009310114507
$ws.Range("K7").Value = 1.035438506210329
$ws.Range("L7").Value = 1.030875884622336
$ws.Range("M7").Value = 1.030517919622258
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027249092137705
$ws.Range("D8").Value = 1.032050250415554
$ws.Range("E8").Value = 1.027360042462999
$ws.Range("F8").Value = 1.025988336264216
$ws.Range("I8").Value = 1.033332201452239
$ws.Range("J8").Value = 1.03228825222211
$ws.Range("K8").Value = 1.034792851025071
$ws.Range("L8").Value = 1.030115923943467
$ws.Range("M8").Value = 1.028748126807683
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.025145477444585
$ws.Range("D9").Value = 1.030469387394885
$ws.Range("E9").Value = 1.025578244392315
$ws.Range("F9").Value = 1.02241873670177
$ws.Range("I9").Value = 1.032850753914902
$ws.Range("J9").Value = 1.031013596911444
$ws.Range("K9").Value = 1.03365011292574
$ws.Range("L9").Value = 1.028775198331817
$ws.Range("M9").Value = 1.025626264256628
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023743507290554
$ws.Range("D10").Value = 1.029415645914166
$ws.Range("E10").Value = 1.024392749465418
$ws.Range("F10").Value = 1.020038966569256
$ws.Range("I10").Value = 1.032523846072668
$ws.Range("J10").Value = 1.030161104586122
$ws.Range("K10").Value = 1.032884935747183
$ws.Range("L10").Value = 1.02788033852597
$ws.Range("M10").Value = 1.023542578400687
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023136543639724
$ws.Range("D11").Value = 1.028959411754222
$ws.Range("E11").Value = 1.023879986478674
$ws.Range("F11").Value = 1.01900842471415
$ws.Range("I11").Value = 1.032380883515426
$ws.Range("J11").Value = 1.029791321157447
$ws.Range("K11").Value = 1.032552813695864
$ws.Range("L11").Value = 1.027492610245735
$ws.Range("M11").Value = 1.022639683371035
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022911105082051
$ws.Range("D12").Value = 1.028789953083317
$ws.Range("E12").Value = 1.023689608743582
$ws.Range("F12").Value = 1.018625616935724
$ws.Range("I12").Value = 1.032327569209608
$ws.Range("J12").Value = 1.029653869660597
$ws.Range("K12").Value = 1.032429329474591
$ws.Range("L12").Value = 1.027348553616519
$ws.Range("M12").Value = 1.022304206083051
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022959461750657
$ws.Range("D13").Value = 1.028826302217151
$ws.Range("E13").Value = 1.02373044154724
$ws.Range("F13").Value = 1.018707731440495
$ws.Range("I13").Value = 1.032339014885282
$ws.Range("J13").Value = 1.029683357884179
$ws.Range("K13").Value = 1.032455822655152
$ws.Range("L13").Value = 1.027379455927575
$ws.Range("M13").Value = 1.022376171789295
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023117908517319
$ws.Range("D14").Value = 1.028945404101343
$ws.Range("E14").Value = 1.023864248053298
$ws.Range("F14").Value = 1.018976782117668
$ws.Range("I14").Value = 1.032376480855313
$ws.Range("J14").Value = 1.029779961362254
$ws.Range("K14").Value = 1.032542608884354
$ws.Range("L14").Value = 1.027480703235221
$ws.Range("M14").Value = 1.022611954820461
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023215534754688
$ws.Range("D15").Value = 1.029018787648911
$ws.Range("E15").Value = 1.023946701972058
$ws.Range("F15").Value = 1.019142550475455
$ws.Range("I15").Value = 1.032399536833285
$ws.Range("J15").Value = 1.029839469037496
$ws.Range("K15").Value = 1.032596064940849
$ws.Range("L15").Value = 1.027543080134836
$ws.Range("M15").Value = 1.02275721488114
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023783791543231
$ws.Range("D16").Value = 1.029445925624629
$ws.Range("E16").Value = 1.024426791785848
$ws.Range("F16").Value = 1.020107357940374
$ws.Range("I16").Value = 1.032533304341645
$ws.Range("J16").Value = 1.030185632227179
$ws.Range("K16").Value = 1.032906960853551
$ws.Range("L16").Value = 1.02790606556446
$ws.Range("M16").Value = 1.023602486554544
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024140270392747
$ws.Range("D17").Value = 1.029713869871128
$ws.Range("E17").Value = 1.02472809084307
$ws.Range("F17").Value = 1.020712529163461
$ws.Range("I17").Value = 1.03261683591544
$ws.Range("J17").Value = 1.030402597583721
$ws.Range("K17").Value = 1.033101764853489
$ws.Range("L17").Value = 1.02813369045166
$ws.Range("M17").Value = 1.024132526908258
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024348207999564
$ws.Range("D18").Value = 1.029870161313785
$ws.Range("E18").Value = 1.024903887882821
$ws.Range("F18").Value = 1.02106550742036
$ws.Range("I18").Value = 1.032665422486093
$ws.Range("J18").Value = 1.030529087223868
$ws.Range("K18").Value = 1.033215313985754
$ws.Range("L18").Value = 1.028266436248843
$ws.Range("M18").Value = 1.024441628547573
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024419110967212
$ws.Range("D19").Value = 1.029923453334217
$ws.Range("E19").Value = 1.02496383934137
$ws.Range("F19").Value = 1.021185862723212
$ws.Range("I19").Value = 1.03268196618339
$ws.Range("J19").Value = 1.030572206307935
$ws.Range("K19").Value = 1.033254018282684
$ws.Range("L19").Value = 1.028311695067336
$ws.Range("M19").Value = 1.024547013863559
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024102022600946
$ws.Range("D20").Value = 1.0296851215613
$ws.Range("E20").Value = 1.024695758675335
$ws.Range("F20").Value = 1.02064760087455
$ws.Range("I20").Value = 1.032607887832814
$ws.Range("J20").Value = 1.030379325727438
$ws.Range("K20").Value = 1.03308087216463
$ws.Range("L20").Value = 1.028109270931128
$ws.Range("M20").Value = 1.024075665035429
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023071249486551
$ws.Range("D21").Value = 1.02891033136692
$ws.Range("E21").Value = 1.023824843014373
$ws.Range("F21").Value = 1.018897553960424
$ws.Range("I21").Value = 1.032365453902756
$ws.Range("J21").Value = 1.029751516736468
$ws.Range("K21").Value = 1.032517055794896
$ws.Range("L21").Value = 1.027450889448804
$ws.Range("M21").Value = 1.022542525446318
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022423246926217
$ws.Range("D22").Value = 1.028423230855556
$ws.Range("E22").Value = 1.023277757360961
$ws.Range("F22").Value = 1.017797119237759
$ws.Range("I22").Value = 1.032211801588057
$ws.Range("J22").Value = 1.029356224380942
$ws.Range("K22").Value = 1.032161871896425
$ws.Range("L22").Value = 1.027036724422622
$ws.Range("M22").Value = 1.021577987883888
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022766757268667
$ws.Range("D23").Value = 1.028681447993198
$ws.Range("E23").Value = 1.023567730860868
$ws.Range("F23").Value = 1.018380492939019
$ws.Range("I23").Value = 1.032293371625675
$ws.Range("J23").Value = 1.029565829798656
$ws.Range("K23").Value = 1.032350226962886
$ws.Range("L23").Value = 1.027256301405668
$ws.Range("M23").Value = 1.022089365244868
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024119305102129
$ws.Range("D24").Value = 1.029698111672382
$ws.Range("E24").Value = 1.024710368022109
$ws.Range("F24").Value = 1.020676939196494
$ws.Range("I24").Value = 1.032611931506309
$ws.Range("J24").Value = 1.030389841470825
$ws.Range("K24").Value = 1.033090312908472
$ws.Range("L24").Value = 1.028120305135259
$ws.Range("M24").Value = 1.024101358658349
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025689235302129
$ws.Range("D25").Value = 1.030878051815847
$ws.Range("E25").Value = 1.026038466769391
$ws.Range("F25").Value = 1.023341548515344
$ws.Range("I25").Value = 1.032976267809461
$ws.Range("J25").Value = 1.031343606153812
$ws.Range("K25").Value = 1.033946130151284
$ws.Range("L25").Value = 1.029121993134438
$ws.Range("M25").Value = 1.02643375289706
